# Applies the Sprint Backlog.xlsx changes described in the commit:
#  - Martin's two Scrum-board cards move from "IN PROGRESS" (col D) to
#    "REVIEWING" (col F), and get annotated with the reviewer's name
#    (Inacio).
#  - Francisco's "Reviews Done" counter (M4) goes from 1 to 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")

# --- Reviews Done counter ---------------------------------------------
$ws.Range("M4").Value = 3

# --- Row 13: "Martin's code smells" moves from D13 to F13 -------------
$ws.Range("D13").Value = $null
$ws.Range("F13").Value = "Martin's code smells (Inacio)"

# --- Row 14: "Martin's design patterns" moves from D14 to F14 ---------
$ws.Range("D14").Value = $null
$ws.Range("F14").Value = "Martin's design patterns (Inacio"
